$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: swap domain/is5FU column order (A<->B) ---
# Writing directly to the header cells keeps the Table1 ListObject column
# names ("tableColumn/@name") in sync automatically.
$ws.Range("A1").Value = "is5FU"
$ws.Range("B1").Value = "domain"

# --- Data rows: swap A (domain) <-> B (is5FU) and refresh the regression
#     stats (columns C:H) for the is5FU = TRUE rows, which were recomputed. ---
$ws.Range("A2").Value = $false
$ws.Range("B2").Value = "Attn/Wkg Mem/Concen"
$ws.Range("C2").Value = -0.141278138722423
$ws.Range("D2").Value = 0.122607944294637
$ws.Range("E2").Value = -1.15227556856283
$ws.Range("F2").Value = 0.249207853153433
$ws.Range("G2").Value = -0.381585293758404
$ws.Range("H2").Value = 0.0990290163135586

$ws.Range("A3").Value = $true
$ws.Range("B3").Value = "Attn/Wkg Mem/Concen"
$ws.Range("C3").Value = 0.121074385563978
$ws.Range("D3").Value = 0.148810734798101
$ws.Range("E3").Value = 0.81361324993285
$ws.Range("F3").Value = 0.41586654492452801
$ws.Range("G3").Value = -0.17058929515324001
$ws.Range("H3").Value = 0.41273806628119702

$ws.Range("A4").Value = $false
$ws.Range("B4").Value = "Exec Fxn"
$ws.Range("C4").Value = 0.192461839673192
$ws.Range("D4").Value = 0.272635872596423
$ws.Range("E4").Value = 0.705929993145433
$ws.Range("F4").Value = 0.480231674336952
$ws.Range("G4").Value = -0.341894651509447
$ws.Range("H4").Value = 0.726818330855831

$ws.Range("A5").Value = $true
$ws.Range("B5").Value = "Exec Fxn"
$ws.Range("C5").Value = 0.026474029820603699
$ws.Range("D5").Value = 0.17395967784973099
$ws.Range("E5").Value = 0.15218486345710799
$ws.Range("F5").Value = 0.87904113211504598
$ws.Range("G5").Value = -0.314480673527058
$ws.Range("H5").Value = 0.36742873316826602

$ws.Range("A6").Value = $false
$ws.Range("B6").Value = "Info Proc Speed"
$ws.Range("C6").Value = 0.157374360349547
$ws.Range("D6").Value = 0.266268438314947
$ws.Range("E6").Value = 0.591036479372001
$ws.Range("F6").Value = 0.554495978923455
$ws.Range("G6").Value = -0.364502188967474
$ws.Range("H6").Value = 0.679250909666569

$ws.Range("A7").Value = $true
$ws.Range("B7").Value = "Info Proc Speed"
$ws.Range("C7").Value = 0.057671398513899699
$ws.Range("D7").Value = 0.20640564412595
$ws.Range("E7").Value = 0.27940804990152401
$ws.Range("F7").Value = 0.77993169406570095
$ws.Range("G7").Value = -0.34687623017875402
$ws.Range("H7").Value = 0.462219027206554

$ws.Range("A8").Value = $false
$ws.Range("B8").Value = "Motor Speed"
$ws.Range("C8").Value = -0.365881131421029
$ws.Range("D8").Value = 0.817869397289788
$ws.Range("E8").Value = -0.447358872496594
$ws.Range("F8").Value = 0.654615965854824
$ws.Range("G8").Value = -1.96887569416649
$ws.Range("H8").Value = 1.23711343132444

$ws.Range("A9").Value = $true
$ws.Range("B9").Value = "Motor Speed"
$ws.Range("C9").Value = -0.080990131792208403
$ws.Range("D9").Value = 0.22781462082584
$ws.Range("E9").Value = -0.35550892870095502
$ws.Range("F9").Value = 0.72220835710266096
$ws.Range("G9").Value = -0.52749858376250203
$ws.Range("H9").Value = 0.365518320178086

$ws.Range("A10").Value = $false
$ws.Range("B10").Value = "Verb Ability/Lang"
$ws.Range("C10").Value = 0.30640565553756
$ws.Range("D10").Value = 0.490161887867891
$ws.Range("E10").Value = 0.625111138016799
$ws.Range("F10").Value = 0.531898118312354
$ws.Range("G10").Value = -0.654293991277667
$ws.Range("H10").Value = 1.26710530235279

$ws.Range("A11").Value = $true
$ws.Range("B11").Value = "Verb Ability/Lang"
$ws.Range("C11").Value = 0.19493031830982999
$ws.Range("D11").Value = 0.26240086608756802
$ws.Range("E11").Value = 0.74287223672797797
$ws.Range("F11").Value = 0.457559035857244
$ws.Range("G11").Value = -0.31936592873392
$ws.Range("H11").Value = 0.70922656535358097

$ws.Range("A12").Value = $false
$ws.Range("B12").Value = "Verb Mem"
$ws.Range("C12").Value = 1.40648071777462
$ws.Range("D12").Value = 0.200751972354477
$ws.Range("E12").Value = 7.00606176506764
$ws.Range("F12").Value = 0.00000000000245119739018773
$ws.Range("G12").Value = 1.01301408213446
$ws.Range("H12").Value = 1.79994735341477

$ws.Range("A13").Value = $true
$ws.Range("B13").Value = "Verb Mem"
$ws.Range("C13").Value = -0.176185460372766
$ws.Range("D13").Value = 0.14869625006656301
$ws.Range("E13").Value = -1.1848682148601399
$ws.Range("F13").Value = 0.23606955639827301
$ws.Range("G13").Value = -0.467624755139391
$ws.Range("H13").Value = 0.115253834393859

$ws.Range("A14").Value = $false
$ws.Range("B14").Value = "Vis Mem"
$ws.Range("C14").Value = 1.01384561099315
$ws.Range("D14").Value = 0.254096875554416
$ws.Range("E14").Value = 3.98999636961704
$ws.Range("F14").Value = 0.0000660743065817391
$ws.Range("G14").Value = 0.515824886322342
$ws.Range("H14").Value = 1.51186633566397

$ws.Range("A15").Value = $true
$ws.Range("B15").Value = "Vis Mem"
$ws.Range("C15").Value = 0.24004338853923801
$ws.Range("D15").Value = 0.17312870543255199
$ws.Range("E15").Value = 1.38650253254943
$ws.Range("F15").Value = 0.16559350749466101
$ws.Range("G15").Value = -0.099282638798608505
$ws.Range("H15").Value = 0.57936941587708501

$ws.Range("A16").Value = $false
$ws.Range("B16").Value = "Visuospatial"
$ws.Range("C16").Value = -0.364503308338097
$ws.Range("D16").Value = 2.18595780227313
$ws.Range("E16").Value = -0.166747641678654
$ws.Range("F16").Value = 0.867568617636766
$ws.Range("G16").Value = -4.64890187251776
$ws.Range("H16").Value = 3.91989525584156

$ws.Range("A17").Value = $true
$ws.Range("B17").Value = "Visuospatial"
$ws.Range("C17").Value = 0.428900104737861
$ws.Range("D17").Value = 0.63320146523395304
$ws.Range("E17").Value = 0.67735172498280904
$ws.Range("F17").Value = 0.49818281988594099
$ws.Range("G17").Value = -0.81215196207867901
$ws.Range("H17").Value = 1.6699521715543999

# --- Highlight the now-reestimated (non-significant) p-values in rows 13 & 15 ---
$ws.Range("F13").Interior.Color = 65535
$ws.Range("F15").Interior.Color = 65535

# --- Column widths follow the swapped content (best effort; engine quantizes to a 
#     7px/MDW grid so it cannot reproduce the exact bestFit float, but this keeps 
#     column A narrow (boolean) and column B wide (domain text), matching the data. ---
$ws.Columns.Item(1).ColumnWidth = 7
$ws.Columns.Item(2).ColumnWidth = 22

# --- Selection moved to F13 (the reestimated Verb Mem p-value) ---
$ws.Range("F13").Select()